# Insert a new data row at row 10 (shifting existing rows 10-81 down to 11-82)
# and populate it with a new "Arveja Verde" price record for
# "Feria Lagunitas de Puerto Montt".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 10:81 down to 11:82, creating a blank row 10.
$ws.Rows("10:10").Insert()

# Fill the new row 10 with the new record's data.
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C10").Value = "Los Lagos"
$ws.Range("D10").Value = 44515
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 100112022
$ws.Range("G10").Value = "Arveja Verde"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 70
$ws.Range("K10").Value = 20000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 20000
$ws.Range("N10").Value = "`$/saco 25 kilos"
$ws.Range("O10").Value = "Región del Maule"
$ws.Range("P10").Value = 800
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
